# Sprint 4 Backlog - Burndown: log Janera's actual time against the
# "Complete Class Diagram" backlog row (row 22).
#
#   - Actual Time (E22)   -> 2.5 hours
#   - Completed By (F22)  -> Janera (matches Assigned Team Member in D22)
#   - Week 2 remaining (I22) -> 0 (task finished, nothing left after week 2)
#
# The dependent rollups (Estimate Totals row 28, and Janera's per-person
# SUMIF summary in row 35) recalculate automatically from these inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E22").Value = 2.5
$ws.Range("F22").Value = "Janera"
$ws.Range("I22").Value = 0

# Reflect where the author was working when they saved.
$null = $ws.Range("D14").Select()
